# Ca va être long
# Insert a new material row ("Seal, O-ring, Copper") into the Materials sheet,
# just above the existing "Seal, O-ring, Elastomer" row, and leave the
# workbook focused on that sheet/cell as the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")

# Insert a new row before row 31 (shifts rows 31:36 down to 32:37).
$ws.Rows.Item(31).Insert()

$ws.Cells.Item(31, 1).Value = "Seal, O-ring, Copper"
$ws.Cells.Item(31, 2).Value = "Plumbing"
$ws.Cells.Item(31, 3).Value = "Oui"

# Make Materials the active sheet and leave the selection where the author
# left it after typing the new row.
$ws.Activate()
$ws.Range("A32").Select()
